# Swap the order of "System" and the email address in the
# "Recorded By" column (G) on the Session Analysis Results sheet.
# Cells that currently read "System, dnasr281@gmail.com" become
# "dnasr281@gmail.com, System". Cells with any other value (e.g. just
# "System", just the email, or empty) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
